$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IssueCreation")

# Replace "ABC News" client references with "HSBC"
$ws.Range("Z3").Value = "HSBC"
$ws.Range("AB3").Value = "Master Service Agreement - HSBC"
